$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Text corrections (date formats, typo fix, trim trailing space) ---
$ws.Range("A2").Value = "09/Jun/2020"
$ws.Range("A3").Value = "10/Jun/2020"
$ws.Range("A4").Value = "11/Jun/2020"
$ws.Range("D3").Value = "378989"
$ws.Range("W2").Value = "Erradicacion"
$ws.Range("W3").Value = "Erradicacion"
$ws.Range("W4").Value = "Erradicacion"

# --- New observation value for row 3 ---
$ws.Range("K3").Value = "NA"

# --- Alignment fix: D2:E3 should be right-aligned like the rest of column D/E ---
$ws.Range("D2:E2").HorizontalAlignment = -4152
$ws.Range("D3:E3").HorizontalAlignment = -4152
